$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14405.5
$ws.Range("C2").Value = 14405.5
$ws.Range("D2").Value = 0.573147927110687
$ws.Range("F2").Value = 14405.5
$ws.Range("G2").Value = 14405.5
$ws.Range("J2").Value = 997
$ws.Range("K2").Value = 314.5
$ws.Range("L2").Value = 170.5
$ws.Range("Q2").Value = 0.5

$ws.Range("B3").Value = 4266
$ws.Range("C3").Value = 4266
$ws.Range("D3").Value = 0.169730245882072
$ws.Range("F3").Value = 4266
$ws.Range("G3").Value = 4266
$ws.Range("K3").Value = 921
$ws.Range("L3").Value = 458
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0

$ws.Range("B4").Value = 1259
$ws.Range("C4").Value = 1259
$ws.Range("D4").Value = 0.0500915095090316
$ws.Range("F4").Value = 1259
$ws.Range("G4").Value = 1259
$ws.Range("K4").Value = 28

$ws.Range("B5").Value = 2962
$ws.Range("C5").Value = 2962
$ws.Range("D5").Value = 0.117848332935466
$ws.Range("F5").Value = 2962
$ws.Range("G5").Value = 2962
$ws.Range("K5").Value = 616
$ws.Range("L5").Value = 325
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 1.5

$ws.Range("B6").Value = 2241.5
$ws.Range("C6").Value = 2241.5
$ws.Range("D6").Value = 0.0891819845627437
$ws.Range("F6").Value = 2241.5
$ws.Range("G6").Value = 2241.5
$ws.Range("J6").Value = 162.5
$ws.Range("K6").Value = 43
$ws.Range("L6").Value = 28.5
